$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in A17:B17 but keep their formatting (style stays as s="1")
$ws.Range("A17:B17").ClearContents()

# Update the selected/active cell shown in the sheet view to L12
$ws.Range("L12").Select()
